$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# --- Simple single-value cell replacements (rows 1-4) ---
$table.Rows.Item(1).Cells(1).Range.Text = "0M"
$table.Rows.Item(2).Cells(1).Range.Text = "0M"
$table.Rows.Item(3).Cells(1).Range.Text = "0M"
$table.Rows.Item(4).Cells(1).Range.Text = "622"

# --- Row 6: 0.00047 -> 0.00089 ---
$table.Rows.Item(6).Cells(1).Range.Text = "0.00089"

# --- Delete rows 7, 8, 9 (0.00013 / 0.00004 / 0.00015) ---
$table.Rows.Item(9).Delete()
$table.Rows.Item(8).Delete()
$table.Rows.Item(7).Delete()

# After deletion: old row10 (0.00018) is now row7 (unchanged value),
# old row11 (0.00021) is now row8, old row12 (0.03922) is now row9.
$table.Rows.Item(8).Cells(1).Range.Text = "0.00006"
$table.Rows.Item(9).Cells(1).Range.Text = "0.00026"

# --- Insert three new rows after row 9 ---
$newVals = @("0.00038", "0.00048", "0.12812")
$afterRow = $table.Rows.Item(9)
foreach ($val in $newVals) {
    if ($afterRow.Index -eq $table.Rows.Count) {
        $newRow = $table.Rows.Add()
    } else {
        $refRow = $table.Rows.Item($afterRow.Index + 1)
        $newRow = $table.Rows.Add($refRow)
    }
    $newRow.Cells(1).Range.Text = $val
    $afterRow = $newRow
}

# --- Collapse the three tab-separated summary rows to single values ---
$table.Rows.Item(44).Cells(1).Range.Text = "99.95"
$table.Rows.Item(45).Cells(1).Range.Text = "0.13"
$table.Rows.Item(46).Cells(1).Range.Text = "240"

Write-Output ("FinalRowCount=" + $table.Rows.Count)
